$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.929.99'
$ws.Range('E2').Value = '  -0.36%  '
$ws.Range('D3').Value = '3.145.24'
$ws.Range('E3').Value = '  +0.34%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''591.28'
$ws.Range('E5').Value = '  +0.46%  '
$ws.Range('D6').Value = '''145.25'
$ws.Range('E6').Value = '  -1.77%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '3.134.02'
$ws.Range('E8').Value = '  +0.28%  '
$ws.Range('E9').Value = '  -0.70%  '
$ws.Range('E10').Value = '  -0.78%  '
$ws.Range('D11').Value = '''5.88'
$ws.Range('E11').Value = '  +2.07%  '
$ws.Range('E12').Value = '  -1.59%  '
$ws.Range('E13').Value = '  -3.19%  '
$ws.Range('D14').Value = '''37.21'
$ws.Range('E14').Value = '  -0.66%  '
$ws.Range('D15').Value = '3.659.04'
$ws.Range('E15').Value = '  +0.17%  '
$ws.Range('E16').Value = '  -1.43%  '
$ws.Range('D17').Value = '''7.32'
$ws.Range('E17').Value = '  +2.00%  '
$ws.Range('D18').Value = '63.759.65'
$ws.Range('E18').Value = '  -0.36%  '
$ws.Range('D19').Value = '3.136.63'
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('D20').Value = '''468.77'
$ws.Range('E20').Value = '  +0.15%  '
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('D23').Value = '''7.53'
$ws.Range('E23').Value = '  -0.57%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').Value = '''12.97'
$ws.Range('E24').Value = '  -2.12%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = '''81.60'
$ws.Range('E25').Value = '  -1.06%  '
$ws.Range('E26').Value = '  +6.26%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').Value = '''9.76'
$ws.Range('E28').Value = '  +8.35%  '
$ws.Range('E29').Value = '  +8.11%  '
$ws.Range('E30').Value = '  -0.34%  '
$ws.Range('D31').Value = '''2.24'
$ws.Range('E31').Value = '  -0.21%  '
$ws.Range('E32').Value = '  -0.04%  '
$ws.Range('D33').Value = '''27.84'
$ws.Range('E33').Value = '  +2.60%  '
$ws.Range('E34').Value = '  +0.57%  '
$ws.Range('D35').Value = '0.0₃0845'
$ws.Range('E35').Value = '  -5.33%  '
$ws.Range('E36').Value = '  +0.67%  '
$ws.Range('D37').Value = '''6.15'
$ws.Range('E37').Value = '  +0.82%  '
$ws.Range('E38').Value = '  -3.83%  '
$ws.Range('E39').Value = '  -6.66%  '
$ws.Range('D40').Value = '''51.61'
$ws.Range('E40').Value = '  +1.25%  '
$ws.Range('D41').Value = '''9.33'
$ws.Range('E41').Value = '  +7.04%  '
$ws.Range('D42').Value = '''453.44'
$ws.Range('E42').Value = '  -0.60%  '
$ws.Range('E43').Value = '  +4.95%  '
$ws.Range('D44').Value = '''0.0372'
$ws.Range('E44').Value = '  -0.61%  '
$ws.Range('D45').Value = '2.910.76'
$ws.Range('E45').Value = '  +0.43%  '
$ws.Range('D46').Value = '''39.66'
$ws.Range('E46').Value = '  +9.98%  '
$ws.Range('D47').Value = '''0.108'
$ws.Range('E47').Value = '  -3.77%  '
$ws.Range('D48').Value = '''132.20'
$ws.Range('E48').Value = '  +5.33%  '
$ws.Range('E50').Value = '  +2.23%  '
$ws.Range('E51').Value = '  -1.11%  '
